$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "27.207.70"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.905.36"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5275"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3827"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07302"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9047"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08073"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.369"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").Value = "1.795.56"
$ws.Range("E15").Value = "  -5.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008678"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "27.250.78"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.127"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.490"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.344"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.742"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.838"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.870"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8094"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.229"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.372"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.701"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01993"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.982"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.621"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1520"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4918"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.636"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "
